$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 136; existing rows 136-172 shift down to 137-173.
$ws.Rows("136:136").Insert()

# Populate the newly inserted row 136 with the new weekly record.
$ws.Cells.Item(136, 1).Value = 9
$ws.Cells.Item(136, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(136, 3).Value = "Metropolitana"
$ws.Cells.Item(136, 4).Value = 45120
$ws.Cells.Item(136, 5).Value = 13
$ws.Cells.Item(136, 6).Value = 100112022
$ws.Cells.Item(136, 7).Value = "Arveja Verde"
$ws.Cells.Item(136, 8).Value = "Perfection"
$ws.Cells.Item(136, 9).Value = "Primera"
$ws.Cells.Item(136, 10).Value = 52
$ws.Cells.Item(136, 11).Value = 26000
$ws.Cells.Item(136, 12).Value = 28000
$ws.Cells.Item(136, 13).Value = 27000
$ws.Cells.Item(136, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(136, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(136, 16).Value = 1080
$ws.Cells.Item(136, 17).Value = 25
$ws.Cells.Item(136, 18).Value = "Hortaliza"
